$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$Bvals = @(17.50317888919334, 16.84879125733272, 16.43507335361597, 16.26372241909267, 16.23511070378862, 16.43277328385421, 17.28015647762706, 18.83800791761071, 19.90849999955314, 20.37758168278804, 20.55251141460298, 20.51495917433761, 20.39202807111112, 20.31637390677404, 19.87747300273864, 19.60354211103058, 19.44430687535155, 19.39010867105873, 19.63287711771608, 20.42821016875099, 20.93221310807184, 20.66470044453176, 19.61962019812038, 18.42888426474833)
$Cvals = @(11.07942455786873, 10.85001874445005, 10.70643155393468, 10.64729319014796, 10.63743738191201, 10.70563644253129, 11.00092820694904, 11.55577663214021, 11.94546786629334, 12.11823920890677, 12.18296948771392, 12.16906026678991, 12.1235787846558, 12.09562827868108, 11.93408216381562, 11.83379132265263, 11.77568561968884, 11.75594120355313, 11.8445113852684, 12.13695700583499, 12.3240196593571, 12.22456769586513, 11.83966623417953, 11.40861025147835)
$Dvals = @(4.944133667519678, 4.913707733208515, 4.894950157265217, 4.887290831353241, 4.886018177948229, 4.894846918752858, 4.933660746441594, 5.009013080883248, 5.06369840805197, 5.088384418048396, 5.097701144788984, 5.095696076470507, 5.089151561001841, 5.085138660448373, 5.062080947922857, 5.047884169919157, 5.039700849704854, 5.036927196375851, 5.049397302446822, 5.091074725201718, 5.118128921544525, 5.103707748270598, 5.048713281583237, 4.988731648774444)
$Evals = @(11.70524948033952, 11.76323444131394, 11.80116870401735, 11.81721351895127, 11.81991316183899, 11.80138271599663, 11.72475888665873, 11.59299367616743, 11.5074551766894, 11.470988053739, 11.45753055120785, 11.46041321420489, 11.46987384465645, 11.47571457739535, 11.50988760467895, 11.53147798287333, 11.54412629373975, 11.54844830661446, 11.52915583442147, 11.4670854781151, 11.42856964173136, 11.44893853586347, 11.53020494351307, 11.62666038224048)
$Fvals = @(24.80316302425874, 24.85419479756033, 24.89489621898645, 24.9138261725084, 24.91711066439516, 24.89514204168989, 24.81880784429946, 24.74390582622751, 24.73499735289434, 24.74103849321756, 24.74478123313568, 24.74391041298667, 24.74131723411847, 24.73991840994905, 24.73480601242019, 24.73425810902024, 24.73489292762911, 24.73527088215525, 24.73421808502912, 24.74203940287781, 24.75563286880457, 24.7476009661322, 24.73423322151087, 24.75610118112085)
$Lvals = @(9.752315673379817, 9.719777783250368, 9.701404343723132, 9.694325670249009, 9.693175083761375, 9.701307217238222, 9.740766662241732, 9.830641276665938, 9.90394426481234, 9.938789744124014, 9.952193066900399, 9.949297289707328, 9.93988831727953, 9.934151921195255, 9.901696568549392, 9.882165072762501, 9.871072692562556, 9.867341551237772, 9.88422963039606, 9.942646374502933, 9.982034329774541, 9.960904160279069, 9.883295818029573, 9.80502428543344)
$Ovals = @(22.16467789400161, 22.26503839062524, 22.33383275546878, 22.36366128648249, 22.36872238322292, 22.33422778110378, 22.19778796857987, 21.98754856570948, 21.86856891161256, 21.82225644401422, 21.80585097232767, 21.80933370131454, 21.82088401945842, 21.82810660143348, 21.87175348980992, 21.9005365882751, 21.91782666785988, 21.92380672344629, 21.89739646469363, 21.81746061996481, 21.77182016282527, 21.79557245513082, 21.89881380297144, 21.82810660143348)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Range("B" + $row).Value = $Bvals[$i]
    $ws.Range("C" + $row).Value = $Cvals[$i]
    $ws.Range("D" + $row).Value = $Dvals[$i]
    $ws.Range("E" + $row).Value = $Evals[$i]
    $ws.Range("F" + $row).Value = $Fvals[$i]
    $ws.Range("L" + $row).Value = $Lvals[$i]
    $ws.Range("O" + $row).Value = $Ovals[$i]
}
